$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 157.61539
$ws.Range("I9").Value = 78.333336
$ws.Range("J9").Value = 225.57143
$ws.Range("K9").Value = 78.333336
$ws.Range("L9").Value = 225.57143
$ws.Range("M9").Value = 90.666664
$ws.Range("N9").Value = -563.57143

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1055.5555
$ws.Range("I12").Value = 2750
$ws.Range("J12").Value = 571.4286
$ws.Range("K12").Value = 2750
$ws.Range("L12").Value = 571.4286
$ws.Range("M12").Value = -2580
$ws.Range("N12").Value = -911.4286

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 9050
$ws.Range("I33").Value = 133.33333
$ws.Range("K33").Value = 133.33333
$ws.Range("M33").Value = 95.66667000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2234.861
$ws.Range("I62").Value = 2270.238
$ws.Range("J62").Value = 2185.3333
$ws.Range("K62").Value = 2270.238
$ws.Range("L62").Value = 2185.3333
$ws.Range("M62").Value = -1646.238
$ws.Range("N62").Value = -3433.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2234.861
$ws.Range("I65").Value = 2270.238
$ws.Range("J65").Value = 2185.3333
$ws.Range("K65").Value = 11351.19
$ws.Range("L65").Value = 10926.6665
$ws.Range("M65").Value = -8231.189999999999
$ws.Range("N65").Value = -17166.6665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 709.36365
$ws.Range("I96").Value = 779.25
$ws.Range("J96").Value = 523
$ws.Range("K96").Value = 2337.75
$ws.Range("L96").Value = 1569
$ws.Range("M96").Value = -964.75
$ws.Range("N96").Value = -4315

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1737.375
$ws.Range("J97").Value = 1842.7142
$ws.Range("L97").Value = 5528.142599999999
$ws.Range("N97").Value = -6520.142599999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1201.8788
$ws.Range("I132").Value = 557.01886
$ws.Range("J132").Value = 3830.923
$ws.Range("K132").Value = 1671.05658
$ws.Range("L132").Value = 11492.769
$ws.Range("M132").Value = 858.9434200000001
$ws.Range("N132").Value = -16552.769

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2660.0154
$ws.Range("I138").Value = 1395.766
$ws.Range("K138").Value = 4187.298000000001
$ws.Range("M138").Value = 952.7019999999993

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 41669716
$ws.Range("I122").Value = 45457644
$ws.Range("J122").Value = 2507
$ws.Range("K122").Value = 136372932
$ws.Range("L122").Value = 7521
$ws.Range("M122").Value = -136370482
$ws.Range("N122").Value = -12421

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1026.8959
$ws.Range("I132").Value = 736.6591
$ws.Range("K132").Value = 2209.9773
$ws.Range("M132").Value = 320.0227

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10418878
$ws.Range("I31").Value = 12501666
$ws.Range("J31").Value = 4937.5
$ws.Range("K31").Value = 12501666
$ws.Range("L31").Value = 4937.5
$ws.Range("M31").Value = -12501371
$ws.Range("N31").Value = -5527.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 10418878
$ws.Range("I34").Value = 12501666
$ws.Range("J34").Value = 4937.5
$ws.Range("K34").Value = 12501666
$ws.Range("L34").Value = 4937.5
$ws.Range("M34").Value = -12501464
$ws.Range("N34").Value = -5341.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 869.04
$ws.Range("I58").Value = 733
$ws.Range("J58").Value = 1073.1
$ws.Range("K58").Value = 733
$ws.Range("L58").Value = 1073.1
$ws.Range("M58").Value = -530
$ws.Range("N58").Value = -1479.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H75").Value = 18000
$ws.Range("J75").Value = 18000
$ws.Range("L75").Value = 18000
$ws.Range("N75").Value = -19996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H78").Value = 18000
$ws.Range("J78").Value = 18000
$ws.Range("L78").Value = 54000
$ws.Range("N78").Value = -63984

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1273.1
$ws.Range("I132").Value = 1144.2142
$ws.Range("J132").Value = 1573.8334
$ws.Range("K132").Value = 3432.6426
$ws.Range("L132").Value = 4721.5002
$ws.Range("M132").Value = -902.6425999999997
$ws.Range("N132").Value = -9781.5002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 896.375
$ws.Range("I134").Value = 897.9048
$ws.Range("J134").Value = 885.6667
$ws.Range("K134").Value = 2693.7144
$ws.Range("L134").Value = 2657.0001
$ws.Range("M134").Value = -158.7143999999998
$ws.Range("N134").Value = -7727.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 869.04
$ws.Range("I136").Value = 733
$ws.Range("J136").Value = 1073.1
$ws.Range("K136").Value = 2199
$ws.Range("L136").Value = 3219.3
$ws.Range("M136").Value = 351
$ws.Range("N136").Value = -8319.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3031.5
$ws.Range("I3").Value = 3030
$ws.Range("J3").Value = 3033
$ws.Range("K3").Value = 9090
$ws.Range("L3").Value = 9099
$ws.Range("M3").Value = -8978
$ws.Range("N3").Value = -9323

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 4450.3228
$ws.Range("I81").Value = 195
$ws.Range("J81").Value = 4743.793
$ws.Range("K81").Value = 585
$ws.Range("L81").Value = 14231.379
$ws.Range("M81").Value = 538
$ws.Range("N81").Value = -16477.379

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 4450.3228
$ws.Range("I84").Value = 195
$ws.Range("J84").Value = 4743.793
$ws.Range("K84").Value = 1755
$ws.Range("L84").Value = 42694.137
$ws.Range("M84").Value = 3861
$ws.Range("N84").Value = -53926.137

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 7464.5
$ws.Range("I87").Value = 7464.5
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 22393.5
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -21145.5
$ws.Range("N87").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 7464.5
$ws.Range("I90").Value = 7464.5
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 67180.5
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -60940.5
$ws.Range("N90").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 481.9
$ws.Range("I97").Value = 482.2857
$ws.Range("J97").Value = 481
$ws.Range("K97").Value = 1446.8571
$ws.Range("L97").Value = 1443
$ws.Range("M97").Value = -950.8571000000002
$ws.Range("N97").Value = -2435

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 10527
$ws.Range("I110").Value = 15000
$ws.Range("K110").Value = 45000
$ws.Range("M110").Value = -40910

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2164.25
$ws.Range("I122").Value = 2248.7
$ws.Range("J122").Value = 2079.8
$ws.Range("K122").Value = 6746.099999999999
$ws.Range("L122").Value = 6239.400000000001
$ws.Range("M122").Value = -4296.099999999999
$ws.Range("N122").Value = -11139.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 27780528
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 33336334
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 100009002
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -100013942

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7353.0586
$ws.Range("I46").Value = 1955.7778
$ws.Range("J46").Value = 13425
$ws.Range("K46").Value = 1955.7778
$ws.Range("L46").Value = 13425
$ws.Range("M46").Value = -1767.7778
$ws.Range("N46").Value = -13801

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2146.2122
$ws.Range("I132").Value = 2242.0193
$ws.Range("J132").Value = 1790.3572
$ws.Range("K132").Value = 6726.0579
$ws.Range("L132").Value = 5371.071599999999
$ws.Range("M132").Value = -4196.0579
$ws.Range("N132").Value = -10431.0716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 80723.08
$ws.Range("I96").Value = 1633.3334
$ws.Range("J96").Value = 104450
$ws.Range("K96").Value = 1633.3334
$ws.Range("L96").Value = 104450
$ws.Range("M96").Value = -260.3334
$ws.Range("N96").Value = -107196

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 16668933
$ws.Range("I132").Value = 23149068
$ws.Range("J132").Value = 5730
$ws.Range("K132").Value = 69447204
$ws.Range("L132").Value = 17190
$ws.Range("M132").Value = -69444674
$ws.Range("N132").Value = -22250
